$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A18").Value = "MicroplasticImages_valid"
$ws.Range("E18").Value = "check_exists_in_zip(MicroplasticImages)"
$ws.Range("B18").Value = "Name of images is correct"
$ws.Range("C18").Value = "mcrplsts_plcy_drft.pdf"
$ws.Range("D18").Value = "error"

$ws.Range("E18").Select()
